# Session 7: Branch and Bound
# Fill in the "Session 7 (Branch&Bound)" test score and the "Test mark" for
# the first student, and add the grader's feedback comment for that session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4").Value = 9
$ws.Range("I4").Value = "NA"
$ws.Range("H5").Value = "Very good. To greatly improve times, you may also comment lines 35 and 37 in Heap.java file. That part is only used to check if a node is already repeated but in this problem that will never happen. The times are going to be much much better without it. Your way of calculating the heuristic value is maybe not so good: whenever you find a solution you will prune all the nodes which score lower than the solution you found (and those nodes may lead to a better solution after processing them)"

$ws.Range("H5:H12").Select()
